$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Informe")

$data = @(
    @{Row=2;  B="4901"; C="6432"},
    @{Row=3;  B="6285"; C="2107"},
    @{Row=4;  B="8826"; C="5895"},
    @{Row=5;  B="3661"; C="140"},
    @{Row=6;  B="1559"; C="1757"},
    @{Row=7;  B="1416"; C="9312"},
    @{Row=8;  B="5554"; C="4505"},
    @{Row=9;  B="8001"; C="2420"},
    @{Row=10; B="2890"; C="3661"},
    @{Row=11; B="3521"; C="7427"},
    @{Row=12; B="3107"; C="2281"},
    @{Row=13; B="4259"; C="9597"},
    @{Row=14; B="6998"; C="900"},
    @{Row=15; B="4946"; C="5127"},
    @{Row=16; B="1909"; C="7235"},
    @{Row=17; B="8551"; C="1639"},
    @{Row=18; B="4694"; C="438"},
    @{Row=19; B="2578"; C="6562"},
    @{Row=20; B="8435"; C="1694"},
    @{Row=21; B="6858"; C="2659"},
    @{Row=22; B="2964"; C="2268"},
    @{Row=23; B="3236"; C="2170"},
    @{Row=24; B="1465"; C="1783"},
    @{Row=25; B="7462"; C="4192"},
    @{Row=26; B="7862"; C="7551"}
)

$lastRow = 26
$ws.Range("B2:D$lastRow").NumberFormat = "@"

foreach ($item in $data) {
    $r = $item.Row
    $b = $item.B
    $c = $item.C
    $d = "$b-$c"

    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
}
